$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("H18").Value = 1622
$ws.Range("I18").Value = 1622
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1622
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -1338
$ws.Range("H138").Value = 8337.903
$ws.Range("I138").Value = 9497.857
$ws.Range("J138").Value = 7999.5835
$ws.Range("K138").Value = 28493.571
$ws.Range("L138").Value = 23998.7505
$ws.Range("M138").Value = -23353.571
$ws.Range("N138").Value = -34278.75049999999

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("H5").Value = 777.5
$ws.Range("I5").Value = 166.25
$ws.Range("J5").Value = 2000
$ws.Range("K5").Value = 166.25
$ws.Range("L5").Value = 2000
$ws.Range("M5").Value = -54.25
$ws.Range("N5").Value = -2224
$ws.Range("H44").Value = 14000
$ws.Range("I44").Value = 14000
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 14000
$ws.Range("N44").Value = 0
$ws.Range("M44").Value = -13512
$ws.Range("H55").Value = 10166.667
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 28500
$ws.Range("K55").Value = 1000
$ws.Range("L55").Value = 28500
$ws.Range("M55").Value = -685
$ws.Range("N55").Value = -29130
$ws.Range("H74").Value = 576323.0600000001
$ws.Range("I74").Value = 651560.75
$ws.Range("J74").Value = 16220.556
$ws.Range("K74").Value = 651560.75
$ws.Range("L74").Value = 16220.556
$ws.Range("M74").Value = -650686.75
$ws.Range("N74").Value = -17968.556
$ws.Range("H77").Value = 576323.0600000001
$ws.Range("I77").Value = 651560.75
$ws.Range("J77").Value = 16220.556
$ws.Range("K77").Value = 3257803.75
$ws.Range("L77").Value = 81102.78
$ws.Range("M77").Value = -3253435.75
$ws.Range("N77").Value = -89838.78
$ws.Range("H80").Value = 30085
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 30085
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 30085
$ws.Range("N80").Value = -32081
$ws.Range("H83").Value = 30085
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 30085
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 90255
$ws.Range("N83").Value = -100239

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("H4").Value = 777.5
$ws.Range("I4").Value = 166.25
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 166.25
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = -51.25
$ws.Range("N4").Value = -2230
$ws.Range("H105").Value = 2352.3809
$ws.Range("I105").Value = 2236
$ws.Range("J105").Value = 2847
$ws.Range("K105").Value = 2236
$ws.Range("L105").Value = 2847
$ws.Range("M105").Value = -489
$ws.Range("N105").Value = -6341

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("H22").Value = 346.47058
$ws.Range("I22").Value = 216.45454
$ws.Range("J22").Value = 584.8333
$ws.Range("K22").Value = 216.45454
$ws.Range("L22").Value = 584.8333
$ws.Range("M22").Value = 133.54546
$ws.Range("N22").Value = -1284.8333
$ws.Range("H41").Value = 20776.166
$ws.Range("I41").Value = 15059
$ws.Range("J41").Value = 21919.6
$ws.Range("K41").Value = 15059
$ws.Range("L41").Value = 21919.6
$ws.Range("M41").Value = -14631
$ws.Range("N41").Value = -22775.6
$ws.Range("H58").Value = 4571002
$ws.Range("I58").Value = 11113954
$ws.Range("J58").Value = 1405057.4
$ws.Range("K58").Value = 11113954
$ws.Range("L58").Value = 1405057.4
$ws.Range("M58").Value = -11113751
$ws.Range("N58").Value = -1405463.4
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("N62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("N65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H132").Value = 3705.8462
$ws.Range("I132").Value = 3878.4546
$ws.Range("J132").Value = 2756.5
$ws.Range("K132").Value = 11635.3638
$ws.Range("L132").Value = 8269.5
$ws.Range("M132").Value = -9105.363799999999
$ws.Range("N132").Value = -13329.5
$ws.Range("H136").Value = 4571002
$ws.Range("I136").Value = 11113954
$ws.Range("J136").Value = 1405057.4
$ws.Range("K136").Value = 33341862
$ws.Range("L136").Value = 4215172.199999999
$ws.Range("M136").Value = -33339312
$ws.Range("N136").Value = -4220272.199999999

# Sheet 5
$ws = $wb.Worksheets.Item(5)
$ws.Range("H68").Value = 5303.6104
$ws.Range("I68").Value = 3625
$ws.Range("J68").Value = 5362.509
$ws.Range("K68").Value = 10875
$ws.Range("L68").Value = 16087.527
$ws.Range("M68").Value = -10064
$ws.Range("N68").Value = -17709.527
$ws.Range("H71").Value = 5303.6104
$ws.Range("I71").Value = 3625
$ws.Range("J71").Value = 5362.509
$ws.Range("K71").Value = 32625
$ws.Range("L71").Value = 48262.581
$ws.Range("M71").Value = -28569
$ws.Range("N71").Value = -56374.581
$ws.Range("H112").Value = 10505.4
$ws.Range("I112").Value = 2351.3333
$ws.Range("J112").Value = 14000
$ws.Range("K112").Value = 7053.999899999999
$ws.Range("L112").Value = 42000
$ws.Range("M112").Value = -5945.999899999999
$ws.Range("N112").Value = -44216
$ws.Range("H131").Value = 3882.8262
$ws.Range("I131").Value = 1537.1111
$ws.Range("J131").Value = 4453.4053
$ws.Range("K131").Value = 4611.3333
$ws.Range("L131").Value = 13360.2159
$ws.Range("M131").Value = 428.6666999999998
$ws.Range("N131").Value = -23440.2159

# Sheet 6
$ws = $wb.Worksheets.Item(6)
$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("N46").Value = 0
$ws.Range("M46").ClearContents()

# Sheet 7
$ws = $wb.Worksheets.Item(7)
$ws.Range("H46").Value = 4200
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -5376
$ws.Range("H47").Value = 28494
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 28494
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 28494
$ws.Range("N47").Value = -29474
$ws.Range("H52").Value = 28494
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 28494
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 28494
$ws.Range("N52").Value = -28960
$ws.Range("H61").Value = 7258.391
$ws.Range("I61").Value = 7392.263
$ws.Range("J61").Value = 6622.5
$ws.Range("K61").Value = 7392.263
$ws.Range("L61").Value = 6622.5
$ws.Range("M61").Value = -7190.263
$ws.Range("N61").Value = -7026.5
$ws.Range("H113").Value = 7258.391
$ws.Range("I113").Value = 7392.263
$ws.Range("J113").Value = 6622.5
$ws.Range("K113").Value = 7392.263
$ws.Range("L113").Value = 6622.5
$ws.Range("M113").Value = -5222.263
$ws.Range("N113").Value = -10962.5

# Sheet 8
$ws = $wb.Worksheets.Item(8)
$ws.Range("H113").Value = 1186.7317
$ws.Range("I113").Value = 741.73914
$ws.Range("J113").Value = 1755.3334
$ws.Range("K113").Value = 2225.21742
$ws.Range("L113").Value = 5266.0002
$ws.Range("M113").Value = -55.21741999999995
$ws.Range("N113").Value = -9606.0002
$ws.Range("H132").Value = 3789805.8
$ws.Range("I132").Value = 3877926
$ws.Range("J132").Value = 635
$ws.Range("K132").Value = 11633778
$ws.Range("L132").Value = 1905
$ws.Range("M132").Value = -11631248
$ws.Range("N132").Value = -6965
$ws.Range("H137").Value = 57957.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 57957.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 57957.5
$ws.Range("N137").Value = -68157.5
